$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.005243333333333
$ws.Range("H2").Value = 3.01573
$ws.Range("I2").Value = 0.07224874268505826
$ws.Range("J2").Value = 0.07224874268505825
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 118.0470123333333
$ws.Range("N2").Value = 354.141037
$ws.Range("O2").Value = 0.4657216250363638
$ws.Range("P2").Value = 0.4657216250363638
$ws.Range("Q2").Value = 118.6659721680011
$ws.Range("R2").Value = 1067.99374951201
$ws.Range("S2").Value = 0.03364780185011944
$ws.Range("T2").Value = 0.03364780185011943

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.005243333333333
$ws.Range("H3").Value = 3.01573
$ws.Range("I3").Value = 0.07224874268505826
$ws.Range("J3").Value = 0.07224874268505825
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("N3").Value = 178.097596
$ws.Range("O3").Value = 0.2342114953037475
$ws.Range("P3").Value = 0.2342114953037476
$ws.Range("Q3").Value = 59.67714035389776
$ws.Range("R3").Value = 537.0942631850799
$ws.Range("S3").Value = 0.01692148605808319
$ws.Range("T3").Value = 0.01692148605808319

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.005243333333333
$ws.Range("H4").Value = 3.01573
$ws.Range("I4").Value = 0.07224874268505826
$ws.Range("J4").Value = 0.07224874268505825
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 56.84506433333333
$ws.Range("N4").Value = 170.535193
$ws.Range("O4").Value = 0.2242663767030476
$ws.Range("P4").Value = 0.2242663767030477
$ws.Range("Q4").Value = 57.14312195398777
$ws.Range("R4").Value = 514.2880975858899
$ws.Range("S4").Value = 0.01620296374332883
$ws.Range("T4").Value = 0.01620296374332883

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.005243333333333
$ws.Range("H5").Value = 3.01573
$ws.Range("I5").Value = 0.07224874268505826
$ws.Range("J5").Value = 0.07224874268505825
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.21324333333333
$ws.Range("N5").Value = 57.63973
$ws.Range("O5").Value = 0.07580050295684103
$ws.Range("P5").Value = 0.07580050295684104
$ws.Range("Q5").Value = 19.31398477254444
$ws.Range("R5").Value = 173.8258629529
$ws.Range("S5").Value = 0.005476491033526805
$ws.Range("T5").Value = 0.005476491033526805

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.25983933333333
$ws.Range("H6").Value = 30.779518
$ws.Range("I6").Value = 0.7373940889775011
$ws.Range("J6").Value = 0.737394088977501
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 118.0470123333333
$ws.Range("N6").Value = 354.141037
$ws.Range("O6").Value = 0.4657216250363638
$ws.Range("P6").Value = 0.4657216250363638
$ws.Range("Q6").Value = 1211.143380320018
$ws.Range("R6").Value = 10900.29042288017
$ws.Range("S6").Value = 0.3434203734108108
$ws.Range("T6").Value = 0.3434203734108108

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 10.25983933333333
$ws.Range("H7").Value = 30.779518
$ws.Range("I7").Value = 0.7373940889775011
$ws.Range("J7").Value = 0.737394088977501
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("N7").Value = 178.097596
$ws.Range("O7").Value = 0.2342114953037475
$ws.Range("P7").Value = 0.2342114953037476
$ws.Range("Q7").Value = 609.084240204303
$ws.Range("R7").Value = 5481.758161838728
$ws.Range("S7").Value = 0.1727061722075652
$ws.Range("T7").Value = 0.1727061722075652

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.25983933333333
$ws.Range("H8").Value = 30.779518
$ws.Range("I8").Value = 0.7373940889775011
$ws.Range("J8").Value = 0.737394088977501
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 56.84506433333333
$ws.Range("N8").Value = 170.535193
$ws.Range("O8").Value = 0.2242663767030476
$ws.Range("P8").Value = 0.2242663767030477
$ws.Range("Q8").Value = 583.2212269529971
$ws.Range("R8").Value = 5248.991042576974
$ws.Range("S8").Value = 0.1653727005372289
$ws.Range("T8").Value = 0.1653727005372289

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.25983933333333
$ws.Range("H9").Value = 30.779518
$ws.Range("I9").Value = 0.7373940889775011
$ws.Range("J9").Value = 0.737394088977501
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.21324333333333
$ws.Range("N9").Value = 57.63973
$ws.Range("O9").Value = 0.07580050295684103
$ws.Range("P9").Value = 0.07580050295684104
$ws.Range("Q9").Value = 197.1247896722378
$ws.Range("R9").Value = 1774.12310705014
$ws.Range("S9").Value = 0.05589484282189617
$ws.Range("T9").Value = 0.05589484282189618

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.383875
$ws.Range("H10").Value = 1.151625
$ws.Range("I10").Value = 0.02758982345723265
$ws.Range("J10").Value = 0.02758982345723265
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 118.0470123333333
$ws.Range("N10").Value = 354.141037
$ws.Range("O10").Value = 0.4657216250363638
$ws.Range("P10").Value = 0.4657216250363638
$ws.Range("Q10").Value = 45.31529685945834
$ws.Range("R10").Value = 407.837671735125
$ws.Range("S10").Value = 0.01284917741496878
$ws.Range("T10").Value = 0.01284917741496878

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.383875
$ws.Range("H11").Value = 1.151625
$ws.Range("I11").Value = 0.02758982345723265
$ws.Range("J11").Value = 0.02758982345723265
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("N11").Value = 178.097596
$ws.Range("O11").Value = 0.2342114953037475
$ws.Range("P11").Value = 0.2342114953037476
$ws.Range("Q11").Value = 22.78907155483333
$ws.Range("R11").Value = 205.1016439935
$ws.Range("S11").Value = 0.006461853807084869
$ws.Range("T11").Value = 0.006461853807084871

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.383875
$ws.Range("H12").Value = 1.151625
$ws.Range("I12").Value = 0.02758982345723265
$ws.Range("J12").Value = 0.02758982345723265
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 56.84506433333333
$ws.Range("N12").Value = 170.535193
$ws.Range("O12").Value = 0.2242663767030476
$ws.Range("P12").Value = 0.2242663767030477
$ws.Range("Q12").Value = 21.82139907095834
$ws.Range("R12").Value = 196.392591638625
$ws.Range("S12").Value = 0.006187469740630318
$ws.Range("T12").Value = 0.006187469740630319

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.383875
$ws.Range("H13").Value = 1.151625
$ws.Range("I13").Value = 0.02758982345723265
$ws.Range("J13").Value = 0.02758982345723265
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.21324333333333
$ws.Range("N13").Value = 57.63973
$ws.Range("O13").Value = 0.07580050295684103
$ws.Range("P13").Value = 0.07580050295684104
$ws.Range("Q13").Value = 7.375483784583334
$ws.Range("R13").Value = 66.37935406125001
$ws.Range("S13").Value = 0.002091322494548686
$ws.Range("T13").Value = 0.002091322494548686

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.264687
$ws.Range("H14").Value = 6.794061
$ws.Range("I14").Value = 0.162767344880208
$ws.Range("J14").Value = 0.162767344880208
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 118.0470123333333
$ws.Range("N14").Value = 354.141037
$ws.Range("O14").Value = 0.4657216250363638
$ws.Range("P14").Value = 0.4657216250363638
$ws.Range("Q14").Value = 267.3395342201396
$ws.Range("R14").Value = 2406.055807981257
$ws.Range("S14").Value = 0.07580427236046473
$ws.Range("T14").Value = 0.07580427236046475

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.264687
$ws.Range("H15").Value = 6.794061
$ws.Range("I15").Value = 0.162767344880208
$ws.Range("J15").Value = 0.162767344880208
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 59.36586533333332
$ws.Range("N15").Value = 178.097596
$ws.Range("O15").Value = 0.2342114953037475
$ws.Range("P15").Value = 0.2342114953037476
$ws.Range("Q15").Value = 134.4451034641506
$ws.Range("R15").Value = 1210.005931177356
$ws.Range("S15").Value = 0.0381219832310143
$ws.Range("T15").Value = 0.0381219832310143

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.264687
$ws.Range("H16").Value = 6.794061
$ws.Range("I16").Value = 0.162767344880208
$ws.Range("J16").Value = 0.162767344880208
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 56.84506433333333
$ws.Range("N16").Value = 170.535193
$ws.Range("O16").Value = 0.2242663767030476
$ws.Range("P16").Value = 0.2242663767030477
$ws.Range("Q16").Value = 128.7362782098637
$ws.Range("R16").Value = 1158.626503888773
$ws.Range("S16").Value = 0.0365032426818596
$ws.Range("T16").Value = 0.0365032426818596

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.264687
$ws.Range("H17").Value = 6.794061
$ws.Range("I17").Value = 0.162767344880208
$ws.Range("J17").Value = 0.162767344880208
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.21324333333333
$ws.Range("N17").Value = 57.63973
$ws.Range("O17").Value = 0.07580050295684103
$ws.Range("P17").Value = 0.07580050295684104
$ws.Range("Q17").Value = 43.51198240483667
$ws.Range("R17").Value = 391.60784164353
$ws.Range("S17").Value = 0.01233784660686937
$ws.Range("T17").Value = 0.01233784660686937

